$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Expand on a few existing to-do items with more specific detail.
$ws.Range("A27").Value = "redo PEA-VFI: figures, interpretations: di/dk, di/dpibar, Obs TR vs anch, PEA vs VFI X1 and X2"
$ws.Range("A28").Value = "redo optimal Taylor rule: figures, table, interpretations: CB loss as a function of psi_pi in RE vs anchoring"
$ws.Range("A26").Value = "get estimation identified: fig alpha hat, autocovariogram"

# 2. New to-do item entered next to the "see loss for RE-optimal TR..." row.
$ws.Range("B20").Value = "Monpol should respond to expectations directly: raise int by x when expectations unanchor, which they do for y forecast errors. "

# 3. A4 - "Susanto's liquidity premium correction..." to-do item is now done:
#    append a note about it and mark it with a new (green-ish) highlight fill
#    instead of the plain yellow one it had before.
$ws.Range("A4").Value = "Susanto's liquidity premium correction for TIPS inflation expectations (write up Notes 27 July and include VIX figure and Andreasen results) - done!"
$ws.Range("A4").Interior.ThemeColor = 10
$ws.Range("A4").WrapText = $true

# 4. Move the active selection to A5 (where the author's cursor ended up).
$ws.Range("A5").Select() | Out-Null
